$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 16.14499125995226
$ws.Range("C2").Value = 16.78066982628543
$ws.Range("D2").Value = 17.44343500719377

$ws.Range("B3").Value = 1.576474232974561
$ws.Range("C3").Value = 1.871035172944628
$ws.Range("D3").Value = 2.478822543852049

$ws.Range("B4").Value = 0.3333225822219811
$ws.Range("C4").Value = 0.3972882749252908
$ws.Range("D4").Value = 0.5251984048125002

$ws.Range("B5").Value = 79.9760060669568
$ws.Range("C5").Value = 80.61050382964908
$ws.Range("D5").Value = 81.61556634214408
